$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.818.83"
$ws.Range("E2").Value = "  -2.51%  "
$ws.Range("D3").Value = "1.782.45"
$ws.Range("E3").Value = "  -2.10%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.84%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5135"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3785"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07770"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -7.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.086"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.198"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "1.771.91"
$ws.Range("E15").Value = "  -2.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.161"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001072"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06533"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.41%  "
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.50%  "
$ws.Range("E22").Value = "  -3.13%  "
$ws.Range("D23").Value = "27.853.26"
$ws.Range("E23").Value = "  -2.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.232"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.38%  "
$ws.Range("D28").Value = "1.979.73"
$ws.Range("E28").Value = "  -2.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.358"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1076"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.027"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.610"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.480"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07092"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02309"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.40%  "
$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2122"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.84%  "
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.642"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.10%  "
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.60%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.020"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6083"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.001"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.150"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.34%  "
$ws.Range("E44").Value = "  -5.71%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5962"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.47%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.712"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.214"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.895"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06704"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.07%  "
